$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 2.971347333333334
$ws.Range("H2").Value = 8.914042
$ws.Range("I2").Value = 0.02922956310646057
$ws.Range("J2").Value = 0.02922956310646057
$ws.Range("M2").Value = 77.08952333333333
$ws.Range("N2").Value = 231.26857
$ws.Range("O2").Value = 0.2403816673726824
$ws.Range("P2").Value = 0.2403816673726824
$ws.Range("Q2").Value = 229.0597495844378
$ws.Range("R2").Value = 2061.53774625994
$ws.Range("S2").Value = 0.007026251116106035
$ws.Range("T2").Value = 0.007026251116106034

$ws.Range("G3").Value = 2.971347333333334
$ws.Range("H3").Value = 8.914042
$ws.Range("I3").Value = 0.02922956310646057
$ws.Range("J3").Value = 0.02922956310646057
$ws.Range("O3").Value = 0.3167483425780597
$ws.Range("P3").Value = 0.3167483425780597
$ws.Range("Q3").Value = 301.8295730503005
$ws.Range("R3").Value = 2716.466157452704
$ws.Range("S3").Value = 0.009258415668252189
$ws.Range("T3").Value = 0.009258415668252187

$ws.Range("G4").Value = 2.971347333333334
$ws.Range("H4").Value = 8.914042
$ws.Range("I4").Value = 0.02922956310646057
$ws.Range("J4").Value = 0.02922956310646057
$ws.Range("O4").Value = 0.4428699900492579
$ws.Range("P4").Value = 0.4428699900492579
$ws.Range("Q4").Value = 422.0109217474952
$ws.Range("R4").Value = 3798.098295727456
$ws.Range("S4").Value = 0.01294489632210235
$ws.Range("T4").Value = 0.01294489632210235

$ws.Range("I5").Value = 0.6284296749927923
$ws.Range("J5").Value = 0.6284296749927923
$ws.Range("M5").Value = 77.08952333333333
$ws.Range("N5").Value = 231.26857
$ws.Range("O5").Value = 0.2403816673726824
$ws.Range("P5").Value = 0.2403816673726824
$ws.Range("Q5").Value = 4924.738131082842
$ws.Range("R5").Value = 44322.64317974559
$ws.Range("S5").Value = 0.1510629731012403
$ws.Range("T5").Value = 0.1510629731012403

$ws.Range("I6").Value = 0.6284296749927923
$ws.Range("J6").Value = 0.6284296749927923
$ws.Range("O6").Value = 0.3167483425780597
$ws.Range("P6").Value = 0.3167483425780597
$ws.Range("S6").Value = 0.1990540579808357
$ws.Range("T6").Value = 0.1990540579808357

$ws.Range("I7").Value = 0.6284296749927923
$ws.Range("J7").Value = 0.6284296749927923
$ws.Range("O7").Value = 0.4428699900492579
$ws.Range("P7").Value = 0.4428699900492579
$ws.Range("S7").Value = 0.2783126439107164
$ws.Range("T7").Value = 0.2783126439107164

$ws.Range("I8").Value = 0.3423407619007471
$ws.Range("J8").Value = 0.3423407619007471
$ws.Range("M8").Value = 77.08952333333333
$ws.Range("N8").Value = 231.26857
$ws.Range("O8").Value = 0.2403816673726824
$ws.Range("P8").Value = 0.2403816673726824
$ws.Range("Q8").Value = 2682.780064413569
$ws.Range("R8").Value = 24145.02057972212
$ws.Range("S8").Value = 0.08229244315533606
$ws.Range("T8").Value = 0.08229244315533606

$ws.Range("I9").Value = 0.3423407619007471
$ws.Range("J9").Value = 0.3423407619007471
$ws.Range("O9").Value = 0.3167483425780597
$ws.Range("P9").Value = 0.3167483425780597
$ws.Range("S9").Value = 0.1084358689289718
$ws.Range("T9").Value = 0.1084358689289718

$ws.Range("I10").Value = 0.3423407619007471
$ws.Range("J10").Value = 0.3423407619007471
$ws.Range("O10").Value = 0.4428699900492579
$ws.Range("P10").Value = 0.4428699900492579
$ws.Range("S10").Value = 0.1516124498164392
$ws.Range("T10").Value = 0.1516124498164392
